$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 16: 2025-12-08 (serial 45999), site 四方坪站
$ws.Range("A16").Value = 45999
$ws.Range("B16").Value = "四方坪站"
$ws.Range("C16").Value = 8839.3700000000008
$ws.Range("D16").Value = 7504.21
$ws.Range("E16").Value = 2976.56
$ws.Range("F16").Value = 362

# New row 17: 2025-12-08 (serial 45999), site 高岭站
$ws.Range("A17").Value = 45999
$ws.Range("B17").Value = "高岭站"
$ws.Range("C17").Value = 4916.45
$ws.Range("D17").Value = 4226.62
$ws.Range("E17").Value = 1287.05
$ws.Range("F17").Value = 185

# Update the selection to mirror the authored view state
$ws.Range("I17").Select()
